$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update footer timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 22:05"

# Row 4
$ws.Range("B4").Value = 1829273
$ws.Range("C4").Value = 12453
$ws.Range("D4").Value = 539173
$ws.Range("E4").Value = 1184166
$ws.Range("G4").Value = 377
$ws.Range("H4").Value = 105934

# Row 5
$ws.Range("B5").Value = 505487
$ws.Range("C5").Value = 7047
$ws.Range("E5").Value = 271103
$ws.Range("G5").Value = 179
$ws.Range("H5").Value = 29013

# Row 9
$ws.Range("B9").Value = 232997
$ws.Range("C9").Value = 333
$ws.Range("E9").Value = 42075

# Row 12
$ws.Range("B12").Value = 183484
$ws.Range("C12").Value = 190
$ws.Range("E12").Value = 9679
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 8605

# Row 17
$ws.Range("B17").Value = 90928
$ws.Range("C17").Value = 738
$ws.Range("D17").Value = 48839
$ws.Range("E17").Value = 34795
$ws.Range("G17").Value = 221
$ws.Range("H17").Value = 7294

# Row 54
$ws.Range("B54").Value = 11398
$ws.Range("C54").Value = 605
$ws.Range("E54").Value = 4706

# Row 72
$ws.Range("B72").Value = 5026
$ws.Range("C72").Value = 226
$ws.Range("D72").Value = 1423
$ws.Range("E72").Value = 3317
$ws.Range("G72").Value = 24
$ws.Range("H72").Value = 286

# Row 84
$ws.Range("B84").Value = 2833
$ws.Range("C84").Value = 34
$ws.Range("D84").Value = 1435
$ws.Range("E84").Value = 1365

# Row 97
$ws.Range("A97").Value = "Maldivas"
$ws.Range("B97").Value = 1773
$ws.Range("C97").Value = 101
$ws.Range("D97").Value = 453
$ws.Range("E97").Value = 1315
$ws.Range("H97").Value = 5

# Row 98
$ws.Range("A98").Value = "Kirguistan"
$ws.Range("B98").Value = 1748
$ws.Range("C98").Value = 26
$ws.Range("D98").Value = 1170
$ws.Range("E98").Value = 562
$ws.Range("H98").Value = 16

# Row 99
$ws.Range("A99").Value = "Mayotte"
$ws.Range("B99").Value = 1699
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 1385
$ws.Range("E99").Value = 293
$ws.Range("H99").Value = 21

# Row 100
$ws.Range("A100").Value = "Lituania"
$ws.Range("B100").Value = 1675
$ws.Range("C100").Value = 5
$ws.Range("D100").Value = 1236
$ws.Range("E100").Value = 369
$ws.Range("H100").Value = 70

# Row 117
$ws.Range("B117").Value = 1056
$ws.Range("C117").Value = 9
$ws.Range("D117").Value = 669
$ws.Range("E117").Value = 377

# Row 143
$ws.Range("A143").Value = "Togo"
$ws.Range("C143").Value = 9
$ws.Range("D143").Value = 211
$ws.Range("E143").Value = 218
$ws.Range("H143").Value = 13

# Row 144
$ws.Range("A144").Value = "Taiwan"
$ws.Range("B144").Value = 442
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 423
$ws.Range("E144").Value = 12
$ws.Range("H144").Value = 7

# Row 145
$ws.Range("A145").Value = "Cabo Verde"
$ws.Range("B145").Value = 435
$ws.Range("C145").Value = 14
$ws.Range("D145").Value = 193
$ws.Range("E145").Value = 238
$ws.Range("H145").Value = 4

# Row 171
$ws.Range("D171").Value = 46
$ws.Range("E171").Value = 71
$ws.Range("G171").Value = 1
$ws.Range("H171").Value = 5

# Row 191
$ws.Range("A191").Value = "Antigua y Barbuda"
$ws.Range("B191").Value = 26
$ws.Range("C191").Value = 1
$ws.Range("D191").Value = 19
$ws.Range("H191").Value = 3

# Row 192
$ws.Range("A192").Value = "Gambia"
$ws.Range("D192").Value = 20
$ws.Range("E192").Value = 4
$ws.Range("H192").Value = 1

# Row 210
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

# Row 211
$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Row 213
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

# Row 214
$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
